# Applies:
#  1. In the "Scenario 1" sample-text paragraph: add explicit sz/szCs=24 paragraph
#     mark run properties, split several runs apart (mirroring Word's automatic
#     proofing-driven run splitting) and wrap "as k" / "kd" with gramStart/gramEnd
#     and spellStart/spellEnd proofErr markers, and fix the trailing
#     "(Figure 2)....\u2026." typo down to a single "(Figure 2)."
#  2. Mark the page setup's orientation explicitly as portrait.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Locate the target paragraph (the "[Sample text] We found that k..." one)
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.StartsWith("[Sample text] We found that k")) {
        $target = $cand
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range

    $newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>[Sample text] We found that k, the numb</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">er of nearest neighbours increased, the naive, brute force performance degraded (see Figure 1).  We hypothesise the reason for this is that </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>as k</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> increases, it takes longer to check each point against the current k-nearest neighbour.  Compare this </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>kd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>-tree p</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>erformance (Figure 2)</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@

    $r.InsertXML($newParaXml)
}

# ---------------------------------------------------------------------------
# 2. Explicitly record the (already-default) portrait page orientation.
# ---------------------------------------------------------------------------
$d.PageSetup.Orientation = 0

Write-Output "edit applied"
